$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new daily rows (06-10-2021 and 07-10-2021) after the last existing row (192)
# The leading apostrophe forces text entry so Excel doesn't reinterpret the
# day/month-ambiguous date strings (e.g. "06-10-2021") as a date serial.
$ws.Cells.Item(193, 1).Value = "'06-10-2021"
$ws.Cells.Item(193, 2).Value = -1.04
$ws.Cells.Item(193, 3).Value = -0.41
$ws.Cells.Item(193, 4).Value = -0.09

$ws.Cells.Item(194, 1).Value = "'07-10-2021"
$ws.Cells.Item(194, 2).Value = -1.13
$ws.Cells.Item(194, 3).Value = -0.45
$ws.Cells.Item(194, 4).Value = -0.15

# Drop the quote-prefix formatting marker so the new cells keep the same
# (default) style as every other data row in the sheet.
$ws.Range("A193:A194").Style = "Normal"
